$wb = $excel.ActiveWorkbook

# The two tabs are mislabelled relative to their actual data: the first
# tab is named "Eintritte" but its content is the "Austrittsinformationen"
# sheet, and the second tab is named "Austritte" but holds the
# "Eintritte (mit Zahlungsinformationen)" data. Fix the tab names (swap
# via a temporary name to avoid a collision while renaming).
$sEintritte = $wb.Worksheets.Item(1)
$sAustritte = $wb.Worksheets.Item(2)

$sEintritte.Name = "TEMP_SWAP_NAME"
$sAustritte.Name = "Eintritte"
$sEintritte.Name = "Austritte"

# Re-fetch sheets by their (now correct) names.
$sAustritte = $wb.Worksheets.Item("Austritte")
$sEintritte = $wb.Worksheets.Item("Eintritte")

# Populate the payment-purpose / address columns (E:K) on row 4 of the
# "Eintritte" sheet (CapitalTransferInformation columns), mirroring the
# header labels already present in row 2.
$sEintritte.Range("E4").Value = "Zusatzname"
$sEintritte.Range("F4").Value = "Str/Postfach"
$sEintritte.Range("G4").Value = "PLZ"
$sEintritte.Range("H4").Value = "Ort"
$sEintritte.Range("I4").Value = "Zahlungszweck"
$sEintritte.Range("J4").Value = "IBAN"
$sEintritte.Range("K4").Value = "ESR Referenznummer"

# Fix up the active tab / selection on each sheet: "Austritte" (formerly
# tabSelected with selection A5) becomes unselected with selection C3;
# "Eintritte" (formerly unselected with selection C3) becomes the active
# tab with selection A5.
$sAustritte.Activate()
$sAustritte.Range("C3").Select()

$sEintritte.Activate()
$sEintritte.Range("A5").Select()
